# Auto-generated Excel COM-interop script applying the Zeromus_Profits diff
# Updates cached market-price / profit values across the ALC, ARM, BSM, CRP, CUL, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 395
$ws.Range("J32").Value = 390
$ws.Range("L32").Value = 390
$ws.Range("N32").Value = -1042
$ws.Range("H40").Value = 2618.1
$ws.Range("J40").Value = 1775
$ws.Range("L40").Value = 1775
$ws.Range("N40").Value = -2125
$ws.Range("H63").Value = 40271
$ws.Range("J63").Value = 40271
$ws.Range("L63").Value = 40271
$ws.Range("N63").Value = -41519
$ws.Range("H64").Value = 2837.3235
$ws.Range("I64").Value = 3042.375
$ws.Range("J64").Value = 2774.2307
$ws.Range("K64").Value = 3042.375
$ws.Range("L64").Value = 2774.2307
$ws.Range("M64").Value = -2794.375
$ws.Range("N64").Value = -3270.2307
$ws.Range("H66").Value = 40271
$ws.Range("J66").Value = 40271
$ws.Range("L66").Value = 120813
$ws.Range("N66").Value = -127053
$ws.Range("H67").Value = 2837.3235
$ws.Range("I67").Value = 3042.375
$ws.Range("J67").Value = 2774.2307
$ws.Range("K67").Value = 3042.375
$ws.Range("L67").Value = 2774.2307
$ws.Range("M67").Value = -2184.375
$ws.Range("N67").Value = -4490.2307
$ws.Range("H70").Value = 1541.6857
$ws.Range("I70").Value = 1700.6296
$ws.Range("J70").Value = 1005.25
$ws.Range("K70").Value = 5101.8888
$ws.Range("L70").Value = 3015.75
$ws.Range("M70").Value = -4831.8888
$ws.Range("N70").Value = -3555.75
$ws.Range("H73").Value = 1541.6857
$ws.Range("I73").Value = 1700.6296
$ws.Range("J73").Value = 1005.25
$ws.Range("K73").Value = 5101.8888
$ws.Range("L73").Value = 3015.75
$ws.Range("M73").Value = -4165.8888
$ws.Range("N73").Value = -4887.75
$ws.Range("H76").Value = 151084.33
$ws.Range("J76").Value = 3866.6667
$ws.Range("L76").Value = 3866.6667
$ws.Range("N76").Value = -4496.6667
$ws.Range("H79").Value = 151084.33
$ws.Range("J79").Value = 3866.6667
$ws.Range("L79").Value = 3866.6667
$ws.Range("N79").Value = -6050.6667
$ws.Range("H127").Value = 482
$ws.Range("I127").Value = 482
$ws.Range("K127").Value = 1446
$ws.Range("M127").Value = 3514
$ws.Range("H129").Value = 939.8823
$ws.Range("I129").Value = 289
$ws.Range("K129").Value = 867
$ws.Range("M129").Value = 4133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2018.3529
$ws.Range("I45").Value = 2119.2727
$ws.Range("J45").Value = 1833.3334
$ws.Range("K45").Value = 2119.2727
$ws.Range("L45").Value = 1833.3334
$ws.Range("M45").Value = -1742.2727
$ws.Range("N45").Value = -2587.3334
$ws.Range("H63").Value = 2720.9583
$ws.Range("I63").Value = 1928.8235
$ws.Range("J63").Value = 4644.7144
$ws.Range("K63").Value = 1928.8235
$ws.Range("L63").Value = 4644.7144
$ws.Range("M63").Value = -1242.8235
$ws.Range("N63").Value = -6016.7144
$ws.Range("H66").Value = 2720.9583
$ws.Range("I66").Value = 1928.8235
$ws.Range("J66").Value = 4644.7144
$ws.Range("K66").Value = 9644.1175
$ws.Range("L66").Value = 23223.572
$ws.Range("M66").Value = -6212.1175
$ws.Range("N66").Value = -30087.572
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = $null
$ws.Range("H132").Value = 2010.6666
$ws.Range("I132").Value = 1265.6154
$ws.Range("J132").Value = 3947.8
$ws.Range("K132").Value = 3796.8462
$ws.Range("L132").Value = 11843.4
$ws.Range("M132").Value = -1266.8462
$ws.Range("N132").Value = -16903.4
$ws.Range("H134").Value = 24389.5
$ws.Range("J134").Value = 24389.5
$ws.Range("L134").Value = 24389.5
$ws.Range("N134").Value = -34529.5
$ws.Range("H135").Value = 30000000
$ws.Range("J135").Value = 30000000
$ws.Range("L135").Value = 30000000
$ws.Range("N135").Value = -30010140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 18.181818
$ws.Range("I22").Value = 18.181818
$ws.Range("K22").Value = 18.181818
$ws.Range("M22").Value = 154.818182
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = $null
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = $null
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").Value = $null
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 52634468
$ws.Range("I62").Value = 2744
$ws.Range("J62").Value = 71431510
$ws.Range("K62").Value = 2744
$ws.Range("L62").Value = 71431510
$ws.Range("M62").Value = -2120
$ws.Range("N62").Value = -71432758
$ws.Range("H65").Value = 52634468
$ws.Range("I65").Value = 2744
$ws.Range("J65").Value = 71431510
$ws.Range("K65").Value = 13720
$ws.Range("L65").Value = 357157550
$ws.Range("M65").Value = -10600
$ws.Range("N65").Value = -357163790
$ws.Range("H99").Value = 1491793.4
$ws.Range("I99").Value = 1704144.8
$ws.Range("J99").Value = 5333.3335
$ws.Range("K99").Value = 1704144.8
$ws.Range("L99").Value = 5333.3335
$ws.Range("M99").Value = -1702646.8
$ws.Range("N99").Value = -8329.333500000001
$ws.Range("H126").Value = 1491793.4
$ws.Range("I126").Value = 1704144.8
$ws.Range("J126").Value = 5333.3335
$ws.Range("K126").Value = 5112434.4
$ws.Range("L126").Value = 16000.0005
$ws.Range("M126").Value = -5109964.4
$ws.Range("N126").Value = -20940.0005
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 4551.5625
$ws.Range("J62").Value = 4688.3335
$ws.Range("L62").Value = 14065.0005
$ws.Range("N62").Value = -15437.0005
$ws.Range("H65").Value = 4551.5625
$ws.Range("J65").Value = 4688.3335
$ws.Range("L65").Value = 42195.0015
$ws.Range("N65").Value = -49059.0015
$ws.Range("H107").Value = 450
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 450
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1350
$ws.Range("M107").Value = $null
$ws.Range("N107").Value = -5190
$ws.Range("H131").Value = 2021343.9
$ws.Range("I131").Value = 66666664
$ws.Range("K131").Value = 199999992
$ws.Range("M131").Value = -199994952
$ws.Range("H132").Value = 719.9583
$ws.Range("I132").Value = 512.9286
$ws.Range("J132").Value = 1009.8
$ws.Range("K132").Value = 4616.3574
$ws.Range("L132").Value = 9088.199999999999
$ws.Range("M132").Value = -2086.3574
$ws.Range("N132").Value = -14148.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 11539.4
$ws.Range("J50").Value = 11539.4
$ws.Range("L50").Value = 11539.4
$ws.Range("N50").Value = -12813.4
$ws.Range("H54").Value = 6833
$ws.Range("J54").Value = 6833
$ws.Range("L54").Value = 6833
$ws.Range("N54").Value = -8121
$ws.Range("H74").Value = 15217
$ws.Range("J74").Value = 15217
$ws.Range("L74").Value = 15217
$ws.Range("N74").Value = -17213
$ws.Range("H77").Value = 15217
$ws.Range("J77").Value = 15217
$ws.Range("L77").Value = 45651
$ws.Range("N77").Value = -55635
$ws.Range("H100").Value = 1892.2222
$ws.Range("I100").Value = 1775
$ws.Range("J100").Value = 1925.7142
$ws.Range("K100").Value = 1775
$ws.Range("L100").Value = 1925.7142
$ws.Range("M100").Value = -1234
$ws.Range("N100").Value = -3007.7142
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1331.8572
$ws.Range("I136").Value = 800.7692
$ws.Range("J136").Value = 2194.875
$ws.Range("K136").Value = 2402.3076
$ws.Range("L136").Value = 6584.625
$ws.Range("M136").Value = 147.6923999999999
$ws.Range("N136").Value = -11684.625
